$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 08:52:39"
$wsZh.Range("H2").Value = "2016-03-20 08:52:59"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 08:52:43"
$wsDe.Range("H2").Value = "2016-03-20 08:53:05"
